$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "Neutrophils"
$ws.Range("G2").Value = 1.519725333333333
$ws.Range("H2").Value = 4.559176
$ws.Range("I2").Value = 0.008503128066097211
$ws.Range("J2").Value = 0.008503128066097211
$ws.Range("M2").Value = 0.264095
$ws.Range("N2").Value = 0.792285
$ws.Range("Q2").Value = 0.4013518619066667
$ws.Range("R2").Value = 3.61216675716
$ws.Range("S2").Value = 0.008503128066097211
$ws.Range("T2").Value = 0.008503128066097211

# Row 3
$ws.Range("D3").Value = "Neutrophils"
$ws.Range("I3").Value = 0.007320459297410611
$ws.Range("J3").Value = 0.007320459297410611
$ws.Range("M3").Value = 0.264095
$ws.Range("N3").Value = 0.792285
$ws.Range("Q3").Value = 0.3455293094716667
$ws.Range("R3").Value = 3.109763785245
$ws.Range("S3").Value = 0.007320459297410611
$ws.Range("T3").Value = 0.007320459297410611

# Row 4
$ws.Range("D4").Value = "Neutrophils"
$ws.Range("G4").Value = 75.66766366666667
$ws.Range("H4").Value = 227.002991
$ws.Range("I4").Value = 0.4233737640003616
$ws.Range("J4").Value = 0.4233737640003616
$ws.Range("M4").Value = 0.264095
$ws.Range("N4").Value = 0.792285
$ws.Range("Q4").Value = 19.98345163604834
$ws.Range("R4").Value = 179.851064724435
$ws.Range("S4").Value = 0.4233737640003616
$ws.Range("T4").Value = 0.4233737640003616

# Row 5
$ws.Range("D5").Value = "Neutrophils"
$ws.Range("G5").Value = 0.7265803333333333
$ws.Range("H5").Value = 2.179741
$ws.Range("I5").Value = 0.00406534357829634
$ws.Range("J5").Value = 0.004065343578296341
$ws.Range("M5").Value = 0.264095
$ws.Range("N5").Value = 0.792285
$ws.Range("Q5").Value = 0.1918862331316667
$ws.Range("R5").Value = 1.726976098185
$ws.Range("S5").Value = 0.00406534357829634
$ws.Range("T5").Value = 0.004065343578296341

# Row 6
$ws.Range("D6").Value = "Neutrophils"
$ws.Range("G6").Value = 44.42260233333334
$ws.Range("H6").Value = 133.267807
$ws.Range("I6").Value = 0.2485522011014548
$ws.Range("J6").Value = 0.2485522011014548
$ws.Range("M6").Value = 0.264095
$ws.Range("N6").Value = 0.792285
$ws.Range("Q6").Value = 11.73178716322167
$ws.Range("R6").Value = 105.586084468995
$ws.Range("S6").Value = 0.2485522011014548
$ws.Range("T6").Value = 0.2485522011014548

# Row 7
$ws.Range("D7").Value = "Neutrophils"
$ws.Range("G7").Value = 55.08051933333333
$ws.Range("H7").Value = 165.241558
$ws.Range("I7").Value = 0.3081851039563794
$ws.Range("J7").Value = 0.3081851039563794
$ws.Range("M7").Value = 0.264095
$ws.Range("N7").Value = 0.792285
$ws.Range("Q7").Value = 14.54648975333667
$ws.Range("R7").Value = 130.91840778003
$ws.Range("S7").Value = 0.3081851039563794
$ws.Range("T7").Value = 0.3081851039563794
